# Generate Report for Handoff
# Adds two new localization rows (cf9a8188... and d3f3f6e5...) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview") -> new rows 6 and 7
# Columns: File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A6").Value = "cf9a8188-ae45-43d9-88db-fcdb7690618d.md"
$wsOverview.Range("B6").Value = "e2e\cf9a8188-ae45-43d9-88db-fcdb7690618d.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-19 04:36:46"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/cf9a8188-ae45-43d9-88db-fcdb7690618d.md", "", "", "e2e\cf9a8188-ae45-43d9-88db-fcdb7690618d.md") | Out-Null

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A7").Value = "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md"
$wsOverview.Range("B7").Value = "e2e\d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-19 04:36:46"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md", "", "", "e2e\d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn") -> new rows 6 and 7
# Columns: Source File Name | File Extension | Status | Source Path | Priority |
#          Content Duplicate | Latest Handoff File | Latest Handoff Datetime |
#          Latest Target File | Latest Handback File | Latest Handback DateTime |
#          Reference Tokens | To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$wsZhCn.Range("A6").Value = "cf9a8188-ae45-43d9-88db-fcdb7690618d.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "'False"
$wsZhCn.Range("G6").Value = "cf9a8188-ae45-43d9-88db-fcdb7690618d.2191c03bad2f0b426731bbde122490c05958718d.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-19 04:36:41"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M6").Value = "'True"
$wsZhCn.Range("O6").Value = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/cf9a8188-ae45-43d9-88db-fcdb7690618d.md", "", "", "cf9a8188-ae45-43d9-88db-fcdb7690618d.md") | Out-Null

$loZhCn.ListRows.Add() | Out-Null
$wsZhCn.Range("A7").Value = "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "'False"
$wsZhCn.Range("G7").Value = "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.8256d8f48a777375261ef6e38096fc63130504f8.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-19 04:36:41"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "'True"
$wsZhCn.Range("O7").Value = "'False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md", "", "", "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de-de") -> new rows 6 and 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$wsDeDe.Range("A6").Value = "cf9a8188-ae45-43d9-88db-fcdb7690618d.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "'False"
$wsDeDe.Range("G6").Value = "cf9a8188-ae45-43d9-88db-fcdb7690618d.2191c03bad2f0b426731bbde122490c05958718d.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-19 04:36:46"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M6").Value = "'True"
$wsDeDe.Range("O6").Value = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/cf9a8188-ae45-43d9-88db-fcdb7690618d.md", "", "", "cf9a8188-ae45-43d9-88db-fcdb7690618d.md") | Out-Null

$loDeDe.ListRows.Add() | Out-Null
$wsDeDe.Range("A7").Value = "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "'False"
$wsDeDe.Range("G7").Value = "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.8256d8f48a777375261ef6e38096fc63130504f8.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-19 04:36:46"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "'True"
$wsDeDe.Range("O7").Value = "'False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md", "", "", "d3f3f6e5-f6ac-41ed-b982-3e10086314c7.md") | Out-Null

Write-Host "Report generated for handoff: added rows 6-7 to Overview, zh-cn and de-de sheets."
